# Apply the "robust" sheet updates: add a parallel SUJ/KAIST "False Positive
# SEGMENTS" comparison block (columns I-N) mirroring the existing
# "Total comparisons" block (columns A-F), and fill in the previously
# missing C13 data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header block (columns I/K), mirrors A2/C2 "Total comparisons" ---
$ws.Range("I1").Value = "False Positive SEGMENTS"
$ws.Range("I2").Value = "Total comparisons"
$ws.Range("K2").Value = 48400

# --- New column headers (SUJ / KAIST), mirrors B4:C4 and E4:F4 ---
$ws.Range("J4").Value = "SUJ"
$ws.Range("K4").Value = "KAIST"
$ws.Range("M4").Value = "SUJ"
$ws.Range("N4").Value = "KAIST"

# --- Row 5 (cif) ---
$ws.Range("I5").Value = "cif"
$ws.Range("J5").Value = 1381
$ws.Range("M5").Formula = "=(J5-2)/(K2)"
$ws.Range("N5").Formula = "=(K5-2)/(K2)"

# --- Row 6 (qcif) ---
$ws.Range("I6").Value = "qcif"
$ws.Range("J6").Value = 1429
$ws.Range("M6").Formula = "=(J6-2)/(K2)"
$ws.Range("N6").Formula = "=(K6-2)/(K2)"

# --- Row 7 (brate) ---
$ws.Range("I7").Value = "brate"
$ws.Range("J7").Value = 1436
$ws.Range("M7").Formula = "=(J7-2)/(K2)"
$ws.Range("N7").Formula = "=(K7-2)/(K2)"

# --- Row 8 (gray) ---
$ws.Range("I8").Value = "gray"
$ws.Range("J8").Value = 1444
$ws.Range("M8").Formula = "=(J8-2)/(K2)"
$ws.Range("N8").Formula = "=(K8-2)/(K2)"

# --- Row 9 (fps) ---
$ws.Range("I9").Value = "fps"
$ws.Range("J9").Value = 935
$ws.Range("M9").Formula = "=(J9-2)/(K2)"
$ws.Range("N9").Formula = "=(K9-2)/(K2)"

# --- Row 10 (5fps) ---
$ws.Range("I10").Value = "5fps"
$ws.Range("J10").Value = 675
$ws.Range("M10").Formula = "=(J10-2)/(K2)"
$ws.Range("N10").Formula = "=(K10-2)/(K2)"

# --- Row 11 (1 degree) ---
$ws.Range("I11").Value = "1 degree"
$ws.Range("J11").Value = 1468
$ws.Range("M11").Formula = "=(J11-2)/(K2)"
$ws.Range("N11").Formula = "=(K11-2)/(K2)"

# --- Row 12 (2 degree) ---
$ws.Range("I12").Value = "2 degree"
$ws.Range("J12").Value = 1450
$ws.Range("M12").Formula = "=(J12-2)/(K2)"
$ws.Range("N12").Formula = "=(K12-2)/(K2)"

# --- Row 13 (3 degree) -- also fill in the previously-missing C13 value ---
$ws.Range("C13").Value = 82
$ws.Range("I13").Value = "3 degree"
$ws.Range("J13").Value = 1418
$ws.Range("M13").Formula = "=(J13-2)/(K2)"
$ws.Range("N13").Formula = "=(K13-2)/(K2)"

# --- Move the selection the way the author's workbook shows it after edit ---
[void]$ws.Range("G14").Select()
